$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 266; this pushes the existing rows
# 266..334 down to 268..336 (dimension grows from R334 to R336),
# matching the weekly refresh described in the commit message.
$ws.Range("A266:A267").EntireRow.Insert()

# New "Primera" observation for this week (row 266).
$ws.Cells.Item(266, 1).Value = 11
$ws.Cells.Item(266, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(266, 3).Value = "Bíobío"
$ws.Cells.Item(266, 4).Value = 45135
$ws.Cells.Item(266, 5).Value = 8
$ws.Cells.Item(266, 6).Value = 100112040
$ws.Cells.Item(266, 7).Value = "Cilantro"
$ws.Cells.Item(266, 8).Value = "Sin especificar"
$ws.Cells.Item(266, 9).Value = "Primera"
$ws.Cells.Item(266, 10).Value = 200
$ws.Cells.Item(266, 11).Value = 600
$ws.Cells.Item(266, 12).Value = 700
$ws.Cells.Item(266, 13).Value = 650
$ws.Cells.Item(266, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(266, 15).Value = "Región de Ñuble"
$ws.Cells.Item(266, 16).Value = 650
$ws.Cells.Item(266, 17).Value = 1
$ws.Cells.Item(266, 18).Value = "Hortaliza"

# New "Segunda" observation for this week (row 267).
$ws.Cells.Item(267, 1).Value = 11
$ws.Cells.Item(267, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(267, 3).Value = "Bíobío"
$ws.Cells.Item(267, 4).Value = 45135
$ws.Cells.Item(267, 5).Value = 8
$ws.Cells.Item(267, 6).Value = 100112040
$ws.Cells.Item(267, 7).Value = "Cilantro"
$ws.Cells.Item(267, 8).Value = "Sin especificar"
$ws.Cells.Item(267, 9).Value = "Segunda"
$ws.Cells.Item(267, 10).Value = 100
$ws.Cells.Item(267, 11).Value = 500
$ws.Cells.Item(267, 12).Value = 500
$ws.Cells.Item(267, 13).Value = 500
$ws.Cells.Item(267, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(267, 15).Value = "Región de Ñuble"
$ws.Cells.Item(267, 16).Value = 500
$ws.Cells.Item(267, 17).Value = 1
$ws.Cells.Item(267, 18).Value = "Hortaliza"
